$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2008.7333
$ws.Range("I15").Value = 2008.7333
$ws.Range("K15").Value = 6026.199900000001
$ws.Range("M15").Value = -5857.199900000001
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H88").Value = 76247.5
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 99996.664
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 99996.664
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -100808.664
$ws.Range("H91").Value = 76247.5
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 99996.664
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 99996.664
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -102804.664
$ws.Range("H93").Value = 25601
$ws.Range("J93").Value = 25601
$ws.Range("L93").Value = 25601
$ws.Range("N93").Value = -30593
$ws.Range("H137").Value = 2847.1304
$ws.Range("I137").Value = 1699.1333
$ws.Range("K137").Value = 5097.3999
$ws.Range("M137").Value = -2547.3999
$ws.Range("H138").Value = 3157.125
$ws.Range("I138").Value = 1439.25
$ws.Range("J138").Value = 4875
$ws.Range("K138").Value = 4317.75
$ws.Range("L138").Value = 14625
$ws.Range("M138").Value = 822.25
$ws.Range("N138").Value = -24905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 10250
$ws.Range("J22").Value = 20000
$ws.Range("L22").Value = 20000
$ws.Range("N22").Value = -20598
$ws.Range("H32").Value = 2316.2559
$ws.Range("I32").Value = 1252.6316
$ws.Range("K32").Value = 1252.6316
$ws.Range("M32").Value = -965.6315999999999
$ws.Range("H102").Value = 1136.3334
$ws.Range("I102").Value = 1204.5
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1204.5
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 417.5
$ws.Range("N102").Value = -4244
$ws.Range("H132").Value = 3648.5833
$ws.Range("I132").Value = 3336.476
$ws.Range("K132").Value = 10009.428
$ws.Range("M132").Value = -7479.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3066.6667
$ws.Range("I105").Value = 3066.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3066.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1319.6667
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1316.6
$ws.Range("I31").Value = 1174.5
$ws.Range("J31").Value = 1529.75
$ws.Range("K31").Value = 1174.5
$ws.Range("L31").Value = 1529.75
$ws.Range("M31").Value = -879.5
$ws.Range("N31").Value = -2119.75
$ws.Range("H34").Value = 1316.6
$ws.Range("I34").Value = 1174.5
$ws.Range("J34").Value = 1529.75
$ws.Range("K34").Value = 1174.5
$ws.Range("L34").Value = 1529.75
$ws.Range("M34").Value = -972.5
$ws.Range("N34").Value = -1933.75
$ws.Range("H58").Value = 1694.8948
$ws.Range("I58").Value = 1613.6
$ws.Range("K58").Value = 1613.6
$ws.Range("M58").Value = -1410.6
$ws.Range("H94").Value = 2059.75
$ws.Range("I94").Value = 1898
$ws.Range("K94").Value = 1898
$ws.Range("M94").Value = -1447
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H136").Value = 1694.8948
$ws.Range("I136").Value = 1613.6
$ws.Range("K136").Value = 4840.799999999999
$ws.Range("M136").Value = -2290.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 103.55556
$ws.Range("J12").Value = 103.55556
$ws.Range("L12").Value = 310.66668
$ws.Range("N12").Value = -656.66668
$ws.Range("H81").Value = 4183.3335
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -17246
$ws.Range("H84").Value = 4183.3335
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 45000
$ws.Range("N84").Value = -56232
$ws.Range("H108").Value = 542
$ws.Range("I108").Value = 542
$ws.Range("K108").Value = 1626
$ws.Range("M108").Value = 1254

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 25751.5
$ws.Range("J20").Value = 25751.5
$ws.Range("L20").Value = 25751.5
$ws.Range("N20").Value = -26241.5
$ws.Range("H24").Value = 1605454.5
$ws.Range("J24").Value = 15238.096
$ws.Range("L24").Value = 15238.096
$ws.Range("N24").Value = -15584.096

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3000.2
$ws.Range("I16").Value = 3000.2
$ws.Range("K16").Value = 3000.2
$ws.Range("M16").Value = -2830.2
$ws.Range("H61").Value = 2958.3333
$ws.Range("I61").Value = 2958.3333
$ws.Range("K61").Value = 2958.3333
$ws.Range("M61").Value = -2756.3333
$ws.Range("H68").Value = 44999.75
$ws.Range("J68").Value = 44999.75
$ws.Range("L68").Value = 44999.75
$ws.Range("N68").Value = -46497.75
$ws.Range("H71").Value = 44999.75
$ws.Range("J71").Value = 44999.75
$ws.Range("L71").Value = 224998.75
$ws.Range("N71").Value = -232486.75
$ws.Range("H113").Value = 2958.3333
$ws.Range("I113").Value = 2958.3333
$ws.Range("K113").Value = 2958.3333
$ws.Range("M113").Value = -788.3332999999998
$ws.Range("H136").Value = 2784.2942
$ws.Range("I136").Value = 2256.3845
$ws.Range("K136").Value = 6769.1535
$ws.Range("M136").Value = -4219.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 12000
$ws.Range("J22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("N22").Value = -12586
$ws.Range("H81").Value = 994
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 994
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 2624.25
$ws.Range("J122").Value = 2998.8
$ws.Range("L122").Value = 8996.400000000001
$ws.Range("N122").Value = -13896.4
$ws.Range("H136").Value = 6868.7856
$ws.Range("I136").Value = 3026.5715
$ws.Range("J136").Value = 10711
$ws.Range("K136").Value = 9079.7145
$ws.Range("L136").Value = 32133
$ws.Range("M136").Value = -6529.7145
$ws.Range("N136").Value = -37233
